$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("F2").Value = 64
$ws.Range("H2").Value = 81
$ws.Range("B9").Value = 279
$ws.Range("E9").Value = 317
$ws.Range("F9").Value = 393
$ws.Range("H9").Value = 322
$ws.Range("I9").Value = 385
$ws.Range("B10").Value = 931
$ws.Range("C10").Value = 1122
$ws.Range("D10").Value = 1277
$ws.Range("E10").Value = 1586
$ws.Range("F10").Value = 1616
$ws.Range("G10").Value = 772
$ws.Range("H10").Value = 403
$ws.Range("I10").Value = 633
$ws.Range("J10").Value = 512
$ws.Range("B11").Value = 1315
$ws.Range("C11").Value = 1592
$ws.Range("D11").Value = 1758
$ws.Range("E11").Value = 2064
$ws.Range("F11").Value = 2174
$ws.Range("G11").Value = 1291
$ws.Range("H11").Value = 918
$ws.Range("I11").Value = 1279
$ws.Range("J11").Value = 1066

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("E8").Value = 12
$ws.Range("E9").Value = 19

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("H2").Value = 2
$ws.Range("H7").Value = 7

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("F2").Value = 4
$ws.Range("F9").Value = 154

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("H7").Value = 19
$ws.Range("H9").Value = 47

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("E7").Value = 50
$ws.Range("I7").Value = 66
$ws.Range("B8").Value = 130
$ws.Range("D8").Value = 375
$ws.Range("E8").Value = 459
$ws.Range("F8").Value = 432
$ws.Range("H8").Value = 70
$ws.Range("I8").Value = 150
$ws.Range("B9").Value = 166
$ws.Range("D9").Value = 433
$ws.Range("E9").Value = 520
$ws.Range("F9").Value = 487
$ws.Range("H9").Value = 132
$ws.Range("I9").Value = 254

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("D6").Value = 6
$ws.Range("D7").Value = 10

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("C6").Value = 20
$ws.Range("J6").Value = 6
$ws.Range("C7").Value = 25
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("H2").Value = 5
$ws.Range("H8").Value = 17

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("B2").Value = 6
$ws.Range("D5").Value = 10
$ws.Range("E8").Value = 74
$ws.Range("G8").Value = 76
$ws.Range("H8").Value = 62
$ws.Range("E11").Value = 8
$ws.Range("B14").Value = 6
$ws.Range("H20").Value = 7
$ws.Range("E21").Value = 19
$ws.Range("H22").Value = 2
$ws.Range("B27").Value = 14
$ws.Range("C28").Value = 104
$ws.Range("F28").Value = 87
$ws.Range("H28").Value = 58
$ws.Range("E29").Value = 19
$ws.Range("F32").Value = 154
$ws.Range("H36").Value = 47
$ws.Range("B47").Value = 36
$ws.Range("B53").Value = 166
$ws.Range("D53").Value = 433
$ws.Range("E53").Value = 520
$ws.Range("F53").Value = 487
$ws.Range("H53").Value = 132
$ws.Range("I53").Value = 254
$ws.Range("H54").Value = 5
$ws.Range("B61").Value = 8
$ws.Range("H65").Value = 17
$ws.Range("F67").Value = 14
$ws.Range("C70").Value = 25
$ws.Range("J70").Value = 13
$ws.Range("I74").Value = 33
$ws.Range("C76").Value = 56
$ws.Range("F77").Value = 45
$ws.Range("I77").Value = 63
$ws.Range("J77").Value = 54
$ws.Range("H83").Value = 14
$ws.Range("E89").Value = 15
$ws.Range("G95").Value = 13
$ws.Range("F97").Value = 18
$ws.Range("B99").Value = 1315
$ws.Range("C99").Value = 1592
$ws.Range("D99").Value = 1758
$ws.Range("E99").Value = 2064
$ws.Range("F99").Value = 2174
$ws.Range("G99").Value = 1291
$ws.Range("H99").Value = 918
$ws.Range("I99").Value = 1279
$ws.Range("J99").Value = 1066

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("E5").Value = 11
$ws.Range("E6").Value = 15

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("C8").Value = 65
$ws.Range("F8").Value = 48
$ws.Range("H8").Value = 21
$ws.Range("C9").Value = 104
$ws.Range("F9").Value = 87
$ws.Range("H9").Value = 58

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("B7").Value = 34
$ws.Range("B8").Value = 36

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("E8").Value = 14
$ws.Range("E9").Value = 19

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("C7").Value = 44
$ws.Range("C8").Value = 56

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("H5").Value = 6
$ws.Range("H6").Value = 14

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I5").Value = 9
$ws.Range("I7").Value = 33

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("F5").Value = 13
$ws.Range("F6").Value = 14

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("F5").Value = 3
$ws.Range("F7").Value = 18

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("B5").Value = 13
$ws.Range("B6").Value = 14

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 6

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("H5").Value = 2
$ws.Range("H6").Value = 5

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("B4").Value = 1
$ws.Range("B6").Value = 6

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("F9").Value = 32
$ws.Range("I9").Value = 29
$ws.Range("J9").Value = 27
$ws.Range("F10").Value = 45
$ws.Range("I10").Value = 63
$ws.Range("J10").Value = 54

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("H6").Value = 2
$ws.Range("H7").Value = 2

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("G6").Value = 10
$ws.Range("G7").Value = 13

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("E6").Value = 6
$ws.Range("E7").Value = 8

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("E7").Value = 41
$ws.Range("G7").Value = 52
$ws.Range("H7").Value = 26
$ws.Range("E8").Value = 74
$ws.Range("G8").Value = 76
$ws.Range("H8").Value = 62
